# Insert the "Helix Render File format" directory tree right after the
# "Python-Rust interop beyond the minimal API." paragraph (and before the
# trailing blank paragraph that precedes the section break).
#
# Each array entry becomes its own new paragraph; an empty string produces
# a blank paragraph (mirrors the blank <w:p/> separators in the target).

$d = $word.ActiveDocument

$lines = @(
    '',
    'Helix Render File format:',
    '',
    'Helix/',
    '├── __init__.py',
    '├── Cargo.toml               # Rust package file',
    '├── pyproject.toml           # Python package config (if using pyo3/maturin)',
    '├── README.md',
    '├── src/                     # Rust source code',
    '│   ├── lib.rs               # Main Rust library entry',
    '│   ├── scene.rs             # Scene class & object management',
    '│   ├── camera.rs            # Camera class',
    '│   ├── mesh.rs              # Cube, Sphere, basic MeshObject',
    '│   ├── renderer.rs          # render() function and GPU interface',
    '│   ├── shader.rs            # Shader loading/compiling',
    '│   └── utils.rs             # Utilities: timing, matrices, etc.',
    '├── render/                  # Python bindings',
    '│   ├── __init__.py',
    '│   ├── scene.py             # Python wrapper for Scene',
    '│   ├── camera.py            # Python wrapper for Camera',
    '│   ├── mesh.py              # Python wrapper for MeshObject, Cube, Sphere',
    '│   └── renderer.py          # Python wrapper for render(), show(), shader interface',
    '├── shaders/                 # GLSL shader files',
    '│   └── basic.glsl',
    '├── examples/                # Example scripts for usage',
    '│   └── basic_scene.py',
    '├── tests/                   # Unit tests',
    '│   ├── test_scene.rs',
    '│   ├── test_camera.rs',
    '│   └── test_renderer.rs',
    '└── target/                  # Rust build output (ignored in git)'
)

# Locate the anchor paragraph ("Python-Rust interop beyond the minimal API.")
# by scanning the paragraphs collection for its text.
$anchorText = "Python-Rust interop beyond the minimal API."
$anchor = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $anchorText) {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Could not find anchor paragraph '$anchorText'"
}

# Put the caret at the end of the anchor paragraph and type the new
# paragraphs in, one at a time, the way a user would (Enter, then text).
$sel = $word.Selection
$anchor.Range.Select()
$sel.Collapse(0)

foreach ($line in $lines) {
    [void]$sel.TypeParagraph()
    if ($line -ne '') {
        [void]$sel.TypeText($line)
    }
}

Write-Output "inserted $($lines.Count) paragraphs after '$anchorText'"
